{"js": "const replacements = [\n  [\"956\u00f73=318, 2\", \"139\u00f76=23, 1\"],\n  [\"911\u00f79=101, 2\", \"344\u00f76=57, 2\"],\n  [\"470\u00f78=58, 6\", \"252\u00f72=126, 0\"],\n  [\"211\u00f78=26, 3\", \"393\u00f72=196, 1\"],\n  [\"446\u00f76=74, 2\", \"392\u00f72=196, 0\"],\n  [\"874\u00f75=174, 4\", \"894\u00f77=127, 5\"],\n  [\"960\u00f78=120, 0\", \"521\u00f74=130, 1\"],\n  [\"359\u00f74=89, 3\", \"689\u00f73=229, 2\"],\n  [\"676\u00f75=135, 1\", \"352\u00f72=176, 0\"],\n  [\"151\u00f79=16, 7\", \"159\u00f75=31, 4\"],\n  [\"268\u00f75=53, 3\", \"523\u00f79=58, 1\"],\n  [\"913\u00f79=101, 4\", \"420\u00f79=46, 6\"],\n  [\"959\u00f78=119, 7\", \"519\u00f73=173, 0\"],\n  [\"964\u00f75=192, 4\", \"324\u00f78=40, 4\"],\n  [\"288\u00f79=32, 0\", \"417\u00f76=69, 3\"],\n  [\"445\u00f77=63, 4\", \"136\u00f72=68, 0\"],\n  [\"929\u00f79=103, 2\", \"288\u00f76=48, 0\"],\n  [\"772\u00f76=128, 4\", \"138\u00f78=17, 2\"],\n  [\"234\u00f74=58, 2\", \"658\u00f72=329, 0\"],\n  [\"273\u00f74=68, 1\", \"665\u00f75=133, 0\"],\n  [\"649\u00f75=129, 4\", \"505\u00f74=126, 1\"],\n  [\"302\u00f73=100, 2\", \"114\u00f72=57, 0\"],\n  [\"225\u00f77=32, 1\", \"874\u00f76=145, 4\"],\n  [\"352\u00f78=44, 0\", \"174\u00f73=58, 0\"],\n  [\"256\u00f76=42, 4\", \"124\u00f79=13, 7\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"956\u00f73=318, 2\", \"139\u00f76=23, 1\"),\n    @(\"911\u00f79=101, 2\", \"344\u00f76=57, 2\"),\n    @(\"470\u00f78=58, 6\", \"252\u00f72=126, 0\"),\n    @(\"211\u00f78=26, 3\", \"393\u00f72=196, 1\"),\n    @(\"446\u00f76=74, 2\", \"392\u00f72=196, 0\"),\n    @(\"874\u00f75=174, 4\", \"894\u00f77=127, 5\"),\n    @(\"960\u00f78=120, 0\", \"521\u00f74=130, 1\"),\n    @(\"359\u00f74=89, 3\", \"689\u00f73=229, 2\"),\n    @(\"676\u00f75=135, 1\", \"352\u00f72=176, 0\"),\n    @(\"151\u00f79=16, 7\", \"159\u00f75=31, 4\"),\n    @(\"268\u00f75=53, 3\", \"523\u00f79=58, 1\"),\n    @(\"913\u00f79=101, 4\", \"420\u00f79=46, 6\"),\n    @(\"959\u00f78=119, 7\", \"519\u00f73=173, 0\"),\n    @(\"964\u00f75=192, 4\", \"324\u00f78=40, 4\"),\n    @(\"288\u00f79=32, 0\", \"417\u00f76=69, 3\"),\n    @(\"445\u00f77=63, 4\", \"136\u00f72=68, 0\"),\n    @(\"929\u00f79=103, 2\", \"288\u00f76=48, 0\"),\n    @(\"772\u00f76=128, 4\", \"138\u00f78=17, 2\"),\n    @(\"234\u00f74=58, 2\", \"658\u00f72=329, 0\"),\n    @(\"273\u00f74=68, 1\", \"665\u00f75=133, 0\"),\n    @(\"649\u00f75=129, 4\", \"505\u00f74=126, 1\"),\n    @(\"302\u00f73=100, 2\", \"114\u00f72=57, 0\"),\n    @(\"225\u00f77=32, 1\", \"874\u00f76=145, 4\"),\n    @(\"352\u00f78=44, 0\", \"174\u00f73=58, 0\"),\n    @(\"256\u00f76=42, 4\", \"124\u00f79=13, 7\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, [ref]$replace, 2)\n}\n"}
